{"js": "// Update the worksheet date header and all 25 multiplication problems\n// from the \"2024-12-22 Sunday\" edition to the \"2024-12-23 Monday\" edition.\nconst replacements = [\n  [\"2024-12-22 Sunday\", \"2024-12-23 Monday\"],\n  [\"693\u00d76=\", \"943\u00d73=\"],\n  [\"714\u00d72=\", \"564\u00d74=\"],\n  [\"832\u00d74=\", \"897\u00d75=\"],\n  [\"448\u00d79=\", \"227\u00d73=\"],\n  [\"313\u00d77=\", \"394\u00d74=\"],\n  [\"274\u00d72=\", \"268\u00d75=\"],\n  [\"620\u00d78=\", \"823\u00d74=\"],\n  [\"781\u00d72=\", \"976\u00d77=\"],\n  [\"529\u00d78=\", \"869\u00d74=\"],\n  [\"955\u00d79=\", \"551\u00d79=\"],\n  [\"461\u00d76=\", \"408\u00d77=\"],\n  [\"673\u00d74=\", \"218\u00d78=\"],\n  [\"302\u00d76=\", \"293\u00d73=\"],\n  [\"791\u00d77=\", \"456\u00d74=\"],\n  [\"843\u00d77=\", \"693\u00d79=\"],\n  [\"674\u00d78=\", \"911\u00d74=\"],\n  [\"211\u00d79=\", \"375\u00d75=\"],\n  [\"526\u00d78=\", \"573\u00d76=\"],\n  [\"305\u00d72=\", \"973\u00d76=\"],\n  [\"439\u00d73=\", \"604\u00d77=\"],\n  [\"117\u00d74=\", \"899\u00d74=\"],\n  [\"771\u00d72=\", \"831\u00d75=\"],\n  [\"468\u00d72=\", \"135\u00d78=\"],\n  [\"525\u00d75=\", \"252\u00d76=\"],\n  [\"431\u00d77=\", \"178\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and all 25 multiplication problems\n# from the \"2024-12-22 Sunday\" edition to the \"2024-12-23 Monday\" edition.\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"2024-12-22 Sunday\", \"2024-12-23 Monday\"),\n    @(\"693\u00d76=\", \"943\u00d73=\"),\n    @(\"714\u00d72=\", \"564\u00d74=\"),\n    @(\"832\u00d74=\", \"897\u00d75=\"),\n    @(\"448\u00d79=\", \"227\u00d73=\"),\n    @(\"313\u00d77=\", \"394\u00d74=\"),\n    @(\"274\u00d72=\", \"268\u00d75=\"),\n    @(\"620\u00d78=\", \"823\u00d74=\"),\n    @(\"781\u00d72=\", \"976\u00d77=\"),\n    @(\"529\u00d78=\", \"869\u00d74=\"),\n    @(\"955\u00d79=\", \"551\u00d79=\"),\n    @(\"461\u00d76=\", \"408\u00d77=\"),\n    @(\"673\u00d74=\", \"218\u00d78=\"),\n    @(\"302\u00d76=\", \"293\u00d73=\"),\n    @(\"791\u00d77=\", \"456\u00d74=\"),\n    @(\"843\u00d77=\", \"693\u00d79=\"),\n    @(\"674\u00d78=\", \"911\u00d74=\"),\n    @(\"211\u00d79=\", \"375\u00d75=\"),\n    @(\"526\u00d78=\", \"573\u00d76=\"),\n    @(\"305\u00d72=\", \"973\u00d76=\"),\n    @(\"439\u00d73=\", \"604\u00d77=\"),\n    @(\"117\u00d74=\", \"899\u00d74=\"),\n    @(\"771\u00d72=\", \"831\u00d75=\"),\n    @(\"468\u00d72=\", \"135\u00d78=\"),\n    @(\"525\u00d75=\", \"252\u00d76=\"),\n    @(\"431\u00d77=\", \"178\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, $wdReplaceAll) | Out-Null\n}\n"}
